$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.658.99"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
$ws.Range("D3").Value = "3.105.69"
$ws.Range("E3").Value = "  +1.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.43"
$ws.Range("E5").Value = "  +1.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.48"
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("D8").Value = "3.112.46"
$ws.Range("E8").Value = "  +1.17%  "

# Row 9
$ws.Range("E9").Value = "  +0.57%  "

# Row 10
$ws.Range("E10").Value = "  +0.26%  "

# Row 11
$ws.Range("E11").Value = "  +0.43%  "

# Row 12
$ws.Range("E12").Value = "  +2.03%  "

# Row 13
$ws.Range("D13").Value = "3.647.79"
$ws.Range("E13").Value = "  +1.66%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.131"
$ws.Range("E14").Value = "  +1.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.22"
$ws.Range("E15").Value = "  +2.43%  "

# Row 16
$ws.Range("E16").Value = "  +0.83%  "

# Row 17
$ws.Range("D17").Value = "57.796.01"
$ws.Range("E17").Value = "  +0.49%  "

# Row 18
$ws.Range("D18").Value = "3.116.14"
$ws.Range("E18").Value = "  +1.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").Value = "  -1.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.05"
$ws.Range("E21").Value = "  -0.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "337.01"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("E24").Value = "  +0.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.27"
$ws.Range("E25").Value = "  +1.08%  "

# Row 26
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
$ws.Range("E27").Value = "  +0.22%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("E28").Value = "  +0.91%  "

# Row 29
$ws.Range("E29").Value = "  +3.73%  "

# Row 30
$ws.Range("E30").Value = "  +0.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("E31").Value = "  +0.93%  "

# Row 32
$ws.Range("E32").Value = "  +2.31%  "

# Row 33
$ws.Range("E33").Value = "  +2.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.90"
$ws.Range("E34").Value = "  +0.18%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.40"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  +2.82%  "

# Row 37
$ws.Range("E37").Value = "  +3.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.92"
$ws.Range("E38").Value = "  -0.45%  "

# Row 39
$ws.Range("E39").Value = "  +0.37%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0665"
$ws.Range("E40").Value = "  -1.20%  "

# Row 41
$ws.Range("D41").Value = "3.156.02"
$ws.Range("E41").Value = "  +1.82%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.683"
$ws.Range("E42").Value = "  +3.99%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.90"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.84"
$ws.Range("E44").Value = "  -0.39%  "

# Row 45
$ws.Range("E45").Value = "  +0.33%  "

# Row 46
$ws.Range("E46").Value = "  +5.66%  "

# Row 47
$ws.Range("D47").Value = "2.277.43"
$ws.Range("E47").Value = "  +0.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0257"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.968"
$ws.Range("E49").Value = "  +4.27%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.63"
$ws.Range("E50").Value = "  +3.49%  "

# Row 51
$ws.Range("E51").Value = "  +2.38%  "
